# Adds a new "2022" column (column P) to worksheet "Лист2", mirroring the
# existing yearly columns D:O (2010:2021).
#
# For every data row that already has a value in column O, this script:
#   1. Copies O's formatting onto the new P cell (Copy + PasteSpecial formats)
#      so the new column matches the existing table's number formats/borders.
#   2. Writes the corresponding 2022 figure into that P cell.
#
# Row 4 (the year header row) gets a literal 2022 label in P4, copying O4's
# header formatting the same way.
#
# Finally, the active selection is moved to Q4 (one cell to the right of the
# newly added column), matching the selection recorded after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - paste only the cell formatting, not any value/clipboard data.
$xlPasteFormats = -4122

# "row|value" pairs: row 4 is the 2022 year header, the rest are the 2022
# data points lining up with each existing data row (columns D:O = 2010:2021).
$data = @(
    "4|2022",
    "6|17.605458686851609",
    "7|23.512267196507903",
    "8|11.567974926285048",
    "9|13.890397233831612",
    "10|15.251616282073513",
    "11|12.556147135845718",
    "12|11.189871602516233",
    "13|14.33643342077432",
    "14|8.0816351285783607",
    "15|11.627388629676959",
    "16|16.762096893820274",
    "17|6.4749486469590067",
    "18|11.869705657244412",
    "19|17.166257166257168",
    "20|6.6785253302338123",
    "21|13.680949038807196",
    "22|15.666264437166006",
    "23|11.723106215150278",
    "24|9.9813900091039063",
    "25|12.54103500424182",
    "26|7.4662568682314392",
    "27|10.470645699175261",
    "28|15.512094411740089",
    "29|5.366655556748146",
    "30|21.623623197012517",
    "31|31.61727133367399",
    "32|9.5604505747835411",
    "33|21.036722921447215",
    "34|29.087239037002597",
    "35|13.358066487477489",
    "38|2.8802473574050205",
    "39|1.5264420238765892",
    "40|4.2640026802302557",
    "41|3.2755413401343145",
    "42|0.70773161401733242",
    "43|5.7924767173652896",
    "44|2.4256163391957477",
    "45|0.84422252477793103",
    "46|3.9877486656962642",
    "47|4.2163214360010697",
    "48|2.2992256801257902",
    "49|6.1400375100473346",
    "50|2.3025931739463203",
    "51|1.1793611793611793",
    "52|3.4034792548306934",
    "53|3.2045466217025864",
    "54|0.93776935011204965",
    "55|5.4399564803481573",
    "56|1.6818459355999253",
    "57|0.66393714728339048",
    "58|2.6820534381025558",
    "59|2.2363577498685321",
    "60|1.6739670228496499",
    "61|2.8057251858086794",
    "62|1.5894585710243447",
    "63|1.9002299437915176",
    "64|1.2143314238132321",
    "65|3.8475059027383725",
    "66|3.0618146354739575",
    "67|4.5969054713991353",
    "70|0.50301310533861732",
    "71|0.846149261854242",
    "72|0.15228581000822344",
    "73|0.52548791552956919",
    "74|0.88466451752166553",
    "75|0.17342744662770332",
    "76|0.28985352480955479",
    "77|0.52188301531726655",
    "78|0.060650169820475497",
    "79|0.27861154863443194",
    "80|0.5562642774497879",
    "81|0",
    "82|0.71347957502562043",
    "83|1.375921375921376",
    "84|0.064216589713786659",
    "85|0.36290805758597666",
    "86|0.64816410963626969",
    "87|0.08159934720522237",
    "88|0.87748483596517857",
    "89|1.475415882851979",
    "90|0.28995172303811417",
    "91|0.30878579809900231",
    "92|0.57658864120376829",
    "93|0.037660740749109793",
    "94|0.13973262162851385",
    "95|0.23952478283086356",
    "96|0.019275101965289396",
    "97|0.44287837729362561",
    "98|0.79380379438213722",
    "99|0.10816248167997966",
    "102|38.882628854480011",
    "103|75.610886033533376",
    "104|1.3418391183743461",
    "105|54.80838958973407",
    "106|109.3799209463787",
    "107|1.3180485943705449",
    "108|37.391104700432571",
    "109|74.214834678205392",
    "110|1.0158903444929646",
    "111|32.78329222265149",
    "112|64.26706618803216",
    "113|1.1907951534637253",
    "114|37.619832137714532",
    "115|74.758394758394758",
    "116|1.2201152045619466",
    "117|38.742146977763326",
    "118|76.276504051025682",
    "119|1.7271861825105399",
    "120|35.794068933746239",
    "121|71.262587141750586",
    "122|0.94234309987387099",
    "123|23.804577889814002",
    "124|46.61068176956914",
    "125|0.71555407423308603",
    "126|27.099395307079902",
    "127|48.479816044966789",
    "128|1.2914318316743894",
    "129|37.063384199760293",
    "130|75.127859111166558",
    "131|0.75713737175985762"
)

foreach ($entry in $data) {
    $parts = $entry.Split("|")
    $r = $parts[0]
    $val = [double]$parts[1]

    $srcCell = $ws.Range("O" + $r)
    $dstCell = $ws.Range("P" + $r)

    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial($xlPasteFormats) | Out-Null
    $dstCell.Value = $val
}

# Match the post-edit selection recorded for the sheet.
$ws.Activate() | Out-Null
$ws.Range("Q4").Select() | Out-Null
